$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registerValidUserSheet")

# Update the test data values (A2 and A3). The CONCAT formulas in D2/D3
# reference these cells, so their cached values will update automatically.
$ws.Range("A2").Value = "testuserAvengers989"
$ws.Range("A3").Value = "Qwerty+12345878889"

# Update the active selection on this sheet to match the authored state.
$ws.Activate()
$ws.Range("C11").Select()
